$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }

    # F holds a date encoded as an 8-digit integer yyyyMMdd.
    # Skip rows whose date does not parse (corrupted data left untouched).
    $fStr = [string][int]$fVal
    if ($fStr.Length -ne 8) {
        continue
    }
    $year = [int]$fStr.Substring(0, 4)
    $month = [int]$fStr.Substring(4, 2)
    $day = [int]$fStr.Substring(6, 2)

    if ($month -lt 1 -or $month -gt 12 -or $day -lt 1 -or $day -gt 31) {
        continue
    }

    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        continue
    }

    $newE = $eVal - 1
    if ($newE -le 0) {
        # Remaining days exhausted: start a new cycle of length D and
        # push the start date forward by D days.
        $newE = $dVal
        $newDate = $startDate.AddDays($dVal)
        $newF = [int]($newDate.ToString("yyyyMMdd"))
        $fCell.Value2 = $newF
    }

    $eCell.Value2 = $newE
}
